$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append new log row 27 (mirrors the existing "Demo inplannen" entries)
$ws.Cells.Item(27, 1).Value = "Demo inplannen"
$ws.Cells.Item(27, 2).Value = "klantenservice@testbedrijf123.nl"
$ws.Cells.Item(27, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Cells.Item(27, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item(27, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Cells.Item(27, 6).Value = "2025-08-13 22:28:12"
$ws.Cells.Item(27, 7).Value = "Nee"
$ws.Cells.Item(27, 8).Value = "Ja"
$ws.Cells.Item(27, 9).Value = "Nee"
$ws.Cells.Item(27, 10).Value = "Nee"

# Extend the conditional formatting ranges from row 26 to row 27
$ws.Range("D2:D26").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D27"))
$ws.Range("G2:G26").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G27"))
$ws.Range("H2:H26").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H27"))
$ws.Range("I2:I26").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I27"))
$ws.Range("J2:J26").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J27"))

# Update the Dashboard category count
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(2, 2).Value = 26
